$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header, matching the bold/border style used by A1:J1 (style index 1)
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "intervention_type"

$values = @(
    "PROCEDURE",
    "DRUG",
    "BIOLOGICAL",
    "BIOLOGICAL",
    "DRUG",
    "DRUG",
    "DRUG",
    "OTHER",
    "BIOLOGICAL",
    "BEHAVIORAL",
    "BEHAVIORAL",
    "OTHER",
    "PROCEDURE",
    "PROCEDURE",
    "PROCEDURE",
    "OTHER",
    "DEVICE",
    "DRUG"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("K$row").Value = $values[$i]
}
